$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 183 (new weekly data entries),
# pushing the existing rows 183:274 down to 185:276.
$ws.Rows.Item(183).Insert()
$ws.Rows.Item(183).Insert()

# Populate the first new row (183) with this week's "Primera" quality data.
$ws.Range("A183").Value2 = 3
$ws.Range("B183").Value2 = "Femacal de La Calera"
$ws.Range("C183").Value2 = "Coquimbo"
$ws.Range("D183").Value2 = 44455
$ws.Range("E183").Value2 = 5
$ws.Range("F183").Value2 = 100112037
$ws.Range("G183").Value2 = "Cebollín"
$ws.Range("H183").Value2 = "Sin especificar"
$ws.Range("I183").Value2 = "Primera"
$ws.Range("J183").Value2 = 150
$ws.Range("K183").Value2 = 3500
$ws.Range("L183").Value2 = 3500
$ws.Range("M183").Value2 = 3500
$ws.Range("N183").Value2 = "`$/paquete 36 unidades"
$ws.Range("O183").Value2 = "Provincia de Quillota"
$ws.Range("P183").Value2 = 97
$ws.Range("Q183").Value2 = 36
$ws.Range("R183").Value2 = "Hortaliza"

# Populate the second new row (184) with this week's "Segunda" quality data.
$ws.Range("A184").Value2 = 3
$ws.Range("B184").Value2 = "Femacal de La Calera"
$ws.Range("C184").Value2 = "Coquimbo"
$ws.Range("D184").Value2 = 44455
$ws.Range("E184").Value2 = 5
$ws.Range("F184").Value2 = 100112037
$ws.Range("G184").Value2 = "Cebollín"
$ws.Range("H184").Value2 = "Sin especificar"
$ws.Range("I184").Value2 = "Segunda"
$ws.Range("J184").Value2 = 160
$ws.Range("K184").Value2 = 2500
$ws.Range("L184").Value2 = 2500
$ws.Range("M184").Value2 = 2500
$ws.Range("N184").Value2 = "`$/paquete 36 unidades"
$ws.Range("O184").Value2 = "Provincia de Quillota"
$ws.Range("P184").Value2 = 69
$ws.Range("Q184").Value2 = 36
$ws.Range("R184").Value2 = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D183:D184").NumberFormat = $ws.Range("D182").NumberFormat
